$d = $word.ActiveDocument

$replacements = @(
    @{old = "73÷8=9, 1"; new = "97÷6=16, 1"},
    @{old = "64÷2=32, 0"; new = "98÷2=49, 0"},
    @{old = "23÷8=2, 7"; new = "31÷7=4, 3"},
    @{old = "86÷6=14, 2"; new = "67÷2=33, 1"},
    @{old = "99÷8=12, 3"; new = "32÷7=4, 4"},
    @{old = "56÷8=7, 0"; new = "65÷9=7, 2"},
    @{old = "41÷3=13, 2"; new = "93÷8=11, 5"},
    @{old = "88÷6=14, 4"; new = "68÷9=7, 5"},
    @{old = "83÷4=20, 3"; new = "13÷3=4, 1"},
    @{old = "20÷2=10, 0"; new = "58÷2=29, 0"},
    @{old = "57÷6=9, 3"; new = "58÷3=19, 1"},
    @{old = "69÷6=11, 3"; new = "83÷9=9, 2"},
    @{old = "35÷2=17, 1"; new = "69÷9=7, 6"},
    @{old = "52÷9=5, 7"; new = "75÷2=37, 1"},
    @{old = "49÷5=9, 4"; new = "47÷3=15, 2"},
    @{old = "28÷4=7, 0"; new = "14÷4=3, 2"},
    @{old = "58÷6=9, 4"; new = "90÷6=15, 0"},
    @{old = "53÷2=26, 1"; new = "45÷3=15, 0"},
    @{old = "84÷9=9, 3"; new = "60÷9=6, 6"},
    @{old = "77÷4=19, 1"; new = "87÷4=21, 3"},
    @{old = "23÷9=2, 5"; new = "82÷6=13, 4"},
    @{old = "30÷2=15, 0"; new = "91÷3=30, 1"},
    @{old = "38÷5=7, 3"; new = "46÷9=5, 1"},
    @{old = "64÷8=8, 0"; new = "96÷7=13, 5"},
    @{old = "76÷3=25, 1"; new = "14÷8=1, 6"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
